$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete "Desarquivamentos Pendentes" sheet (no longer needed)
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
